# Insert a new weekly observation row at row 70 (Apio / Terminal Hortofrutícola
# Agro Chillán). This pushes the existing rows 70-166 down to 71-167, which is
# exactly what the target diff shows (every row below the insertion point keeps
# its original data, just one row lower; dimension grows from R166 to R167).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above row 70; everything that was row 70.. shifts to 71..
$ws.Rows(70).Insert()

# The freshly inserted row 70 is blank. The new observation reuses every field
# of the record that is now sitting in row 71 (old row 70) except the date and
# the volume, so clone that row into place first …
$ws.Range("A71:R71").Copy()
$ws.Range("A70").PasteSpecial()

# … then overwrite the two fields that actually changed for the new entry.
$ws.Range("D70").Value = 44571   # 2022-01-10
$ws.Range("J70").Value = 60
